$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 9 (shifts existing rows 9-60 down to 10-61),
# inheriting formatting (e.g. the date style on column D) from the row below.
$ws.Rows(9).Insert()

# Populate the newly inserted row with the latest weekly price report.
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 44635
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 100112012
$ws.Range("G9").Value = "Espinaca"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 170
$ws.Range("K9").Value = 7500
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 7765
$ws.Range("N9").Value = "$/cuna 10 kilos"
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 776
$ws.Range("Q9").Value = 10
$ws.Range("R9").Value = "Hortaliza"
